# Update column G ("K" - strikeouts) values for rows 2 through 45 on Sheet1.
# The underlying source data changed (switch from "Strike#" to "K" metric),
# so these are freshly regenerated literal values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 0
    4  = 1
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 4
    14 = 1
    15 = 1
    16 = 1
    17 = 1
    18 = 0
    19 = 2
    20 = 2
    21 = 1
    22 = 0
    23 = 3
    24 = 0
    25 = 2
    26 = 0
    27 = 0
    28 = 1
    29 = 0
    30 = 2
    31 = 1
    32 = 1
    33 = 2
    34 = 1
    35 = 0
    36 = 1
    37 = 2
    38 = 0
    39 = 0
    40 = 1
    41 = 1
    42 = 2
    43 = 1
    44 = 1
    45 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
